$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 39.288329
$ws.Cells.Item(2, 8).Value2 = 117.864987
$ws.Cells.Item(2, 9).Value2 = 0.632237668435316
$ws.Cells.Item(2, 10).Value2 = 0.632237668435316
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 112.513392
$ws.Cells.Item(2, 14).Value2 = 337.540176
$ws.Cells.Item(2, 15).Value2 = 0.3275312977368564
$ws.Cells.Item(2, 16).Value2 = 0.3275312977368564
$ws.Cells.Item(2, 17).Value2 = 4420.463161801968
$ws.Cells.Item(2, 18).Value2 = 39784.16845621771
$ws.Cells.Item(2, 19).Value2 = 0.2070776240207434
$ws.Cells.Item(2, 20).Value2 = 0.2070776240207434

# Row 3
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 39.288329
$ws.Cells.Item(3, 8).Value2 = 117.864987
$ws.Cells.Item(3, 9).Value2 = 0.632237668435316
$ws.Cells.Item(3, 10).Value2 = 0.632237668435316
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 106.314466
$ws.Cells.Item(3, 14).Value2 = 318.943398
$ws.Cells.Item(3, 15).Value2 = 0.3094859589441663
$ws.Cells.Item(3, 16).Value2 = 0.3094859589441664
$ws.Cells.Item(3, 17).Value2 = 4176.917717667314
$ws.Cells.Item(3, 18).Value2 = 37592.25945900583
$ws.Cells.Item(3, 19).Value2 = 0.1956686810963277
$ws.Cells.Item(3, 20).Value2 = 0.1956686810963277

# Row 4
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 39.288329
$ws.Cells.Item(4, 8).Value2 = 117.864987
$ws.Cells.Item(4, 9).Value2 = 0.632237668435316
$ws.Cells.Item(4, 10).Value2 = 0.632237668435316
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 124.6916553333333
$ws.Cells.Item(4, 14).Value2 = 374.074966
$ws.Cells.Item(4, 15).Value2 = 0.3629827433189773
$ws.Cells.Item(4, 16).Value2 = 0.3629827433189773
$ws.Cells.Item(4, 17).Value2 = 4898.926778290604
$ws.Cells.Item(4, 18).Value2 = 44090.34100461545
$ws.Cells.Item(4, 19).Value2 = 0.229491363318245
$ws.Cells.Item(4, 20).Value2 = 0.229491363318245

# Row 5
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 19.344283
$ws.Cells.Item(5, 8).Value2 = 58.032849
$ws.Cells.Item(5, 9).Value2 = 0.3112930657211948
$ws.Cells.Item(5, 10).Value2 = 0.3112930657211947
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 112.513392
$ws.Cells.Item(5, 14).Value2 = 337.540176
$ws.Cells.Item(5, 15).Value2 = 0.3275312977368564
$ws.Cells.Item(5, 16).Value2 = 0.3275312977368564
$ws.Cells.Item(5, 17).Value2 = 2176.490896137936
$ws.Cells.Item(5, 18).Value2 = 19588.41806524142
$ws.Cells.Item(5, 19).Value2 = 0.1019582217921474
$ws.Cells.Item(5, 20).Value2 = 0.1019582217921474

# Row 6
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 19.344283
$ws.Cells.Item(6, 8).Value2 = 58.032849
$ws.Cells.Item(6, 9).Value2 = 0.3112930657211948
$ws.Cells.Item(6, 10).Value2 = 0.3112930657211947
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 106.314466
$ws.Cells.Item(6, 14).Value2 = 318.943398
$ws.Cells.Item(6, 15).Value2 = 0.3094859589441663
$ws.Cells.Item(6, 16).Value2 = 0.3094859589441664
$ws.Cells.Item(6, 17).Value2 = 2056.577117297878
$ws.Cells.Item(6, 18).Value2 = 18509.1940556809
$ws.Cells.Item(6, 19).Value2 = 0.09634083295739335
$ws.Cells.Item(6, 20).Value2 = 0.09634083295739335

# Row 7
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 19.344283
$ws.Cells.Item(7, 8).Value2 = 58.032849
$ws.Cells.Item(7, 9).Value2 = 0.3112930657211948
$ws.Cells.Item(7, 10).Value2 = 0.3112930657211947
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 124.6916553333333
$ws.Cells.Item(7, 14).Value2 = 374.074966
$ws.Cells.Item(7, 15).Value2 = 0.3629827433189773
$ws.Cells.Item(7, 16).Value2 = 0.3629827433189773
$ws.Cells.Item(7, 17).Value2 = 2412.07066850646
$ws.Cells.Item(7, 18).Value2 = 21708.63601655814
$ws.Cells.Item(7, 19).Value2 = 0.112994010971654
$ws.Cells.Item(7, 20).Value2 = 0.112994010971654

# Row 8
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 3.509096666666667
$ws.Cells.Item(8, 8).Value2 = 10.52729
$ws.Cells.Item(8, 9).Value2 = 0.05646926584348937
$ws.Cells.Item(8, 10).Value2 = 0.05646926584348937
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 112.513392
$ws.Cells.Item(8, 14).Value2 = 337.540176
$ws.Cells.Item(8, 15).Value2 = 0.3275312977368564
$ws.Cells.Item(8, 16).Value2 = 0.3275312977368564
$ws.Cells.Item(8, 17).Value2 = 394.82036882256
$ws.Cells.Item(8, 18).Value2 = 3553.38331940304
$ws.Cells.Item(8, 19).Value2 = 0.01849545192396561
$ws.Cells.Item(8, 20).Value2 = 0.01849545192396561

# Row 9
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 3.509096666666667
$ws.Cells.Item(9, 8).Value2 = 10.52729
$ws.Cells.Item(9, 9).Value2 = 0.05646926584348937
$ws.Cells.Item(9, 10).Value2 = 0.05646926584348937
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 106.314466
$ws.Cells.Item(9, 14).Value2 = 318.943398
$ws.Cells.Item(9, 15).Value2 = 0.3094859589441663
$ws.Cells.Item(9, 16).Value2 = 0.3094859589441664
$ws.Cells.Item(9, 17).Value2 = 373.0677382590467
$ws.Cells.Item(9, 18).Value2 = 3357.60964433142
$ws.Cells.Item(9, 19).Value2 = 0.01747644489044536
$ws.Cells.Item(9, 20).Value2 = 0.01747644489044537

# Row 10
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 3.509096666666667
$ws.Cells.Item(10, 8).Value2 = 10.52729
$ws.Cells.Item(10, 9).Value2 = 0.05646926584348937
$ws.Cells.Item(10, 10).Value2 = 0.05646926584348937
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 124.6916553333333
$ws.Cells.Item(10, 14).Value2 = 374.074966
$ws.Cells.Item(10, 15).Value2 = 0.3629827433189773
$ws.Cells.Item(10, 16).Value2 = 0.3629827433189773
$ws.Cells.Item(10, 17).Value2 = 437.555072091349
$ws.Cells.Item(10, 18).Value2 = 3937.995648822141
$ws.Cells.Item(10, 19).Value2 = 0.02049736902907839
$ws.Cells.Item(10, 20).Value2 = 0.02049736902907839
